# Lab testing report finalized
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TestNo updated
$ws.Range("A2").Value = 200

# Result: Fail -> Pass (ItemType in J2 stays "Fabric", unaffected)
$ws.Range("I2").Value = "Pass"

# New Burst Strength Test data row values
$ws.Range("L6").Value = 200
$ws.Range("M6").Value = 150

# Update the active selection to M6, as recorded in the sheet view
$ws.Range("M6").Select()

$wb.Save()
